$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's default column got a hair narrower in the authored edit
# (15.1173469387755 -> 14.8469387755102 raw sheet units). Nudge it the same
# direction; Excel's COM layer only exposes pixel-quantised character widths
# so this is the closest attainable approximation.
$ws.Columns.Item(1).ColumnWidth = 14

# Add a "Combinations" column header and the COMBIN() helper values used to
# normalise the RTP[%] calculation.
$ws.Range("I10").Value = "Combinations"
$ws.Range("I11").Formula = "=COMBIN(3,3)"
$ws.Range("I12").Formula = "=COMBIN(4,3)"
$ws.Range("I13").Formula = "=COMBIN(5,3)"

# Correct the RTP[%] formulas so they divide by the number of combinations.
$ws.Range("I17").Formula = "=100 * SUM(E17:G17) / I11"
$ws.Range("I18").Formula = "=100 * SUM(E18:G18) / I12"
$ws.Range("I19").Formula = "=100 * SUM(E19:G19) / I13"

# Move the active selection to the newly added cell, as in the authored edit.
$ws.Range("I11").Select() | Out-Null
